$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shared-string text fix: "Trizol" -> "TRIzol".
#    Write the corrected text to every cell that currently shares that string
#    (G2:G27) so the engine updates the single shared-string entry in place
#    instead of forking a new, separate string.
$ws.Range("G2:G27").Value = "TRIzol"

# 2. Bring G3:G27 to the same cell style already used by G2 (Arial 11, black)
#    instead of the old Calibri 12 style.
for ($i = 3; $i -le 27; $i++) {
    $cell = $ws.Range("G$i")
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
}

# 3. H2:H27 become real FALSE() formulas (rather than bare boolean literals).
for ($i = 2; $i -le 27; $i++) {
    $ws.Range("H$i").Formula = "=FALSE()"
}

# 4. Move the sheet's active selection from H2:H27 to G2:G27.
$ws.Range("G2:G27").Select()
